$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 835.0909
$ws.Cells.Item(4, 10).Value = 1999
$ws.Cells.Item(4, 12).Value = 1999
$ws.Cells.Item(4, 14).Value = -2227
$ws.Cells.Item(33, 8).Value = 190.77777
$ws.Cells.Item(33, 9).Value = 195.875
$ws.Cells.Item(33, 11).Value = 195.875
$ws.Cells.Item(33, 13).Value = 33.125
$ws.Cells.Item(61, 8).Value = 124.5
$ws.Cells.Item(61, 9).Value = 124.5
$ws.Cells.Item(61, 11).Value = 373.5
$ws.Cells.Item(61, 13).Value = -201.5
$ws.Cells.Item(62, 8).Value = 3202.2
$ws.Cells.Item(62, 9).Value = 1999.5
$ws.Cells.Item(62, 10).Value = 4004
$ws.Cells.Item(62, 11).Value = 1999.5
$ws.Cells.Item(62, 12).Value = 4004
$ws.Cells.Item(62, 13).Value = -1375.5
$ws.Cells.Item(62, 14).Value = -5252
$ws.Cells.Item(65, 8).Value = 3202.2
$ws.Cells.Item(65, 9).Value = 1999.5
$ws.Cells.Item(65, 10).Value = 4004
$ws.Cells.Item(65, 11).Value = 9997.5
$ws.Cells.Item(65, 12).Value = 20020
$ws.Cells.Item(65, 13).Value = -6877.5
$ws.Cells.Item(65, 14).Value = -26260
$ws.Cells.Item(70, 8).Value = 1261
$ws.Cells.Item(70, 10).Value = 987.5
$ws.Cells.Item(70, 12).Value = 2962.5
$ws.Cells.Item(70, 14).Value = -3502.5
$ws.Cells.Item(73, 8).Value = 1261
$ws.Cells.Item(73, 10).Value = 987.5
$ws.Cells.Item(73, 12).Value = 2962.5
$ws.Cells.Item(73, 14).Value = -4834.5
$ws.Cells.Item(121, 8).Value = 2085.875
$ws.Cells.Item(121, 10).Value = 2085.875
$ws.Cells.Item(121, 12).Value = 6257.625
$ws.Cells.Item(121, 14).Value = -9751.625
$ws.Cells.Item(132, 8).Value = 7674.914
$ws.Cells.Item(132, 9).Value = 5109.8486
$ws.Cells.Item(132, 11).Value = 15329.5458
$ws.Cells.Item(132, 13).Value = -12799.5458
$ws.Cells.Item(135, 8).Value = 2044.8334
$ws.Cells.Item(135, 9).Value = 550
$ws.Cells.Item(135, 10).Value = 5034.5
$ws.Cells.Item(135, 11).Value = 4950
$ws.Cells.Item(135, 12).Value = 45310.5
$ws.Cells.Item(135, 13).Value = -2415
$ws.Cells.Item(135, 14).Value = -50380.5
$ws.Cells.Item(137, 8).Value = 7158772
$ws.Cells.Item(137, 9).Value = 16667467
$ws.Cells.Item(137, 10).Value = 27251.25
$ws.Cells.Item(137, 11).Value = 50002401
$ws.Cells.Item(137, 12).Value = 81753.75
$ws.Cells.Item(137, 13).Value = -49999851
$ws.Cells.Item(137, 14).Value = -86853.75
$ws.Cells.Item(138, 8).Value = 387589.94
$ws.Cells.Item(138, 9).Value = 2867.1333
$ws.Cells.Item(138, 10).Value = 717352.3
$ws.Cells.Item(138, 11).Value = 8601.3999
$ws.Cells.Item(138, 12).Value = 2152056.9
$ws.Cells.Item(138, 13).Value = -3461.3999
$ws.Cells.Item(138, 14).Value = -2162336.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(33, 8).Value = 28000
$ws.Cells.Item(33, 9).Value = 28000
$ws.Cells.Item(33, 11).Value = 28000
$ws.Cells.Item(33, 13).Value = -27671
$ws.Cells.Item(61, 8).Value = 2679.4849
$ws.Cells.Item(61, 9).Value = 1623.1482
$ws.Cells.Item(61, 11).Value = 1623.1482
$ws.Cells.Item(61, 13).Value = -1411.1482
$ws.Cells.Item(74, 8).Value = 200934.17
$ws.Cells.Item(74, 9).Value = 328487.94
$ws.Cells.Item(74, 10).Value = 3805.6365
$ws.Cells.Item(74, 11).Value = 328487.94
$ws.Cells.Item(74, 12).Value = 3805.6365
$ws.Cells.Item(74, 13).Value = -327613.94
$ws.Cells.Item(74, 14).Value = -5553.636500000001
$ws.Cells.Item(77, 8).Value = 200934.17
$ws.Cells.Item(77, 9).Value = 328487.94
$ws.Cells.Item(77, 10).Value = 3805.6365
$ws.Cells.Item(77, 11).Value = 1642439.7
$ws.Cells.Item(77, 12).Value = 19028.1825
$ws.Cells.Item(77, 13).Value = -1638071.7
$ws.Cells.Item(77, 14).Value = -27764.1825
$ws.Cells.Item(122, 8).Value = 3407.7354
$ws.Cells.Item(122, 9).Value = 3226.6333
$ws.Cells.Item(122, 11).Value = 9679.8999
$ws.Cells.Item(122, 13).Value = -7229.8999
$ws.Cells.Item(132, 8).Value = 3106.6155
$ws.Cells.Item(132, 9).Value = 1876.3334
$ws.Cells.Item(132, 11).Value = 5629.0002
$ws.Cells.Item(132, 13).Value = -3099.0002
$ws.Cells.Item(136, 8).Value = 2679.4849
$ws.Cells.Item(136, 9).Value = 1623.1482
$ws.Cells.Item(136, 11).Value = 4869.444600000001
$ws.Cells.Item(136, 13).Value = -2319.444600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 3000
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 14).Value = -1700
$ws.Cells.Item(31, 8).Value = 3678.2354
$ws.Cells.Item(31, 9).Value = 2620.6365
$ws.Cells.Item(31, 11).Value = 2620.6365
$ws.Cells.Item(31, 13).Value = -2325.6365
$ws.Cells.Item(34, 8).Value = 3678.2354
$ws.Cells.Item(34, 9).Value = 2620.6365
$ws.Cells.Item(34, 11).Value = 2620.6365
$ws.Cells.Item(34, 13).Value = -2418.6365
$ws.Cells.Item(58, 8).Value = 2814.0334
$ws.Cells.Item(58, 9).Value = 2277.3333
$ws.Cells.Item(58, 10).Value = 3350.7334
$ws.Cells.Item(58, 11).Value = 2277.3333
$ws.Cells.Item(58, 12).Value = 3350.7334
$ws.Cells.Item(58, 13).Value = -2074.3333
$ws.Cells.Item(58, 14).Value = -3756.7334
$ws.Cells.Item(122, 8).Value = 3895.7856
$ws.Cells.Item(122, 9).Value = 3038.25
$ws.Cells.Item(122, 11).Value = 9114.75
$ws.Cells.Item(122, 13).Value = -6664.75
$ws.Cells.Item(132, 8).Value = 2771.7144
$ws.Cells.Item(132, 9).Value = 2377.4
$ws.Cells.Item(132, 10).Value = 3297.4666
$ws.Cells.Item(132, 11).Value = 7132.200000000001
$ws.Cells.Item(132, 12).Value = 9892.399800000001
$ws.Cells.Item(132, 13).Value = -4602.200000000001
$ws.Cells.Item(132, 14).Value = -14952.3998
$ws.Cells.Item(136, 8).Value = 2814.0334
$ws.Cells.Item(136, 9).Value = 2277.3333
$ws.Cells.Item(136, 10).Value = 3350.7334
$ws.Cells.Item(136, 11).Value = 6831.999899999999
$ws.Cells.Item(136, 12).Value = 10052.2002
$ws.Cells.Item(136, 13).Value = -4281.999899999999
$ws.Cells.Item(136, 14).Value = -15152.2002
$ws.Cells.Item(139, 8).Value = 98997
$ws.Cells.Item(139, 10).Value = 98997
$ws.Cells.Item(139, 12).Value = 98997
$ws.Cells.Item(139, 14).Value = -109277
$ws.Cells.Item(140, 8).Value = 98570.71000000001
$ws.Cells.Item(140, 10).Value = 103499.164
$ws.Cells.Item(140, 12).Value = 103499.164
$ws.Cells.Item(140, 14).Value = -113859.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1056.4
$ws.Cells.Item(5, 9).Value = 559.2727
$ws.Cells.Item(5, 10).Value = 2423.5
$ws.Cells.Item(5, 11).Value = 1677.8181
$ws.Cells.Item(5, 12).Value = 7270.5
$ws.Cells.Item(5, 13).Value = -1565.8181
$ws.Cells.Item(5, 14).Value = -7494.5
$ws.Cells.Item(135, 8).Value = 1056.4
$ws.Cells.Item(135, 9).Value = 559.2727
$ws.Cells.Item(135, 10).Value = 2423.5
$ws.Cells.Item(135, 11).Value = 5033.454299999999
$ws.Cells.Item(135, 12).Value = 21811.5
$ws.Cells.Item(135, 13).Value = -2498.454299999999
$ws.Cells.Item(135, 14).Value = -26881.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 254587.38
$ws.Cells.Item(70, 9).Value = 403139.8
$ws.Cells.Item(70, 10).Value = 7000
$ws.Cells.Item(70, 11).Value = 403139.8
$ws.Cells.Item(70, 12).Value = 7000
$ws.Cells.Item(70, 13).Value = -402869.8
$ws.Cells.Item(70, 14).Value = -7540
$ws.Cells.Item(73, 8).Value = 254587.38
$ws.Cells.Item(73, 9).Value = 403139.8
$ws.Cells.Item(73, 10).Value = 7000
$ws.Cells.Item(73, 11).Value = 403139.8
$ws.Cells.Item(73, 12).Value = 7000
$ws.Cells.Item(73, 13).Value = -402203.8
$ws.Cells.Item(73, 14).Value = -8872
$ws.Cells.Item(97, 8).Value = 3184.923
$ws.Cells.Item(97, 9).Value = 3125.8333
$ws.Cells.Item(97, 10).Value = 3235.5715
$ws.Cells.Item(97, 11).Value = 3125.8333
$ws.Cells.Item(97, 12).Value = 3235.5715
$ws.Cells.Item(97, 13).Value = -2629.8333
$ws.Cells.Item(97, 14).Value = -4227.5715
$ws.Cells.Item(113, 8).Value = 5779.8335
$ws.Cells.Item(113, 9).Value = 4998
$ws.Cells.Item(113, 11).Value = 4998
$ws.Cells.Item(113, 13).Value = -2828
$ws.Cells.Item(122, 8).Value = 5508.3
$ws.Cells.Item(122, 9).Value = 4302.385
$ws.Cells.Item(122, 10).Value = 7747.857
$ws.Cells.Item(122, 11).Value = 12907.155
$ws.Cells.Item(122, 12).Value = 23243.571
$ws.Cells.Item(122, 13).Value = -10457.155
$ws.Cells.Item(122, 14).Value = -28143.571
$ws.Cells.Item(126, 8).Value = 10527.5
$ws.Cells.Item(126, 9).Value = 8603
$ws.Cells.Item(126, 11).Value = 25809
$ws.Cells.Item(126, 13).Value = -23339
$ws.Cells.Item(132, 8).Value = 5477
$ws.Cells.Item(132, 9).Value = 1800
$ws.Cells.Item(132, 11).Value = 5400
$ws.Cells.Item(132, 13).Value = -2870

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 27589
$ws.Cells.Item(56, 9).Value = 17883.5
$ws.Cells.Item(56, 11).Value = 17883.5
$ws.Cells.Item(56, 13).Value = -17192.5
$ws.Cells.Item(61, 8).Value = 2203.2693
$ws.Cells.Item(61, 9).Value = 2134.8262
$ws.Cells.Item(61, 11).Value = 2134.8262
$ws.Cells.Item(61, 13).Value = -1932.8262
$ws.Cells.Item(113, 8).Value = 2203.2693
$ws.Cells.Item(113, 9).Value = 2134.8262
$ws.Cells.Item(113, 11).Value = 2134.8262
$ws.Cells.Item(113, 13).Value = 35.17380000000003
$ws.Cells.Item(132, 8).Value = 5457.909
$ws.Cells.Item(132, 9).Value = 3089.4167
$ws.Cells.Item(132, 10).Value = 8300.1
$ws.Cells.Item(132, 11).Value = 9268.250100000001
$ws.Cells.Item(132, 12).Value = 24900.3
$ws.Cells.Item(132, 13).Value = -6738.250100000001
$ws.Cells.Item(132, 14).Value = -29960.3
$ws.Cells.Item(138, 8).Value = 84125
$ws.Cells.Item(138, 10).Value = 84125
$ws.Cells.Item(138, 12).Value = 84125
$ws.Cells.Item(138, 14).Value = -94405

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 27500700
$ws.Cells.Item(2, 9).Value = 27500700
$ws.Cells.Item(2, 11).Value = 27500700
$ws.Cells.Item(2, 13).Value = -27500588
$ws.Cells.Item(4, 8).Value = 1002
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(56, 8).Value = 34999.6
$ws.Cells.Item(56, 10).Value = 34999.6
$ws.Cells.Item(56, 12).Value = 34999.6
$ws.Cells.Item(56, 14).Value = -36427.6
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
$ws.Cells.Item(100, 8).Value = 55556948
$ws.Cells.Item(100, 9).Value = 1418.75
$ws.Cells.Item(100, 10).Value = 166668020
$ws.Cells.Item(100, 11).Value = 2837.5
$ws.Cells.Item(100, 12).Value = 333336040
$ws.Cells.Item(100, 13).Value = -2296.5
$ws.Cells.Item(100, 14).Value = -333337122
$ws.Cells.Item(107, 8).Value = 713.2727
$ws.Cells.Item(107, 9).Value = 429.6
$ws.Cells.Item(107, 10).Value = 949.6667
$ws.Cells.Item(107, 11).Value = 1288.8
$ws.Cells.Item(107, 12).Value = 2849.0001
$ws.Cells.Item(107, 13).Value = 631.1999999999998
$ws.Cells.Item(107, 14).Value = -6689.0001
$ws.Cells.Item(113, 8).Value = 1090.2
$ws.Cells.Item(113, 9).Value = 1302.5
$ws.Cells.Item(113, 10).Value = 948.6667
$ws.Cells.Item(113, 11).Value = 3907.5
$ws.Cells.Item(113, 12).Value = 2846.0001
$ws.Cells.Item(113, 13).Value = -1737.5
$ws.Cells.Item(113, 14).Value = -7186.0001
